$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 11; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $val = [string]$cell.Value2
    if ($val.EndsWith("%")) {
        $newVal = $val.Substring(0, $val.Length - 1)
        $cell.Formula = "'" + $newVal
    }
}
